# Inserts a new price-record row at row 193 of the "Hortaliza, Terminal La
# Palmera de La Serena - Zanahoria" sheet. All existing rows from 193 to 256
# shift down by one (to 194..257); the new row 193 holds a fresh weekly
# record (date 2021-12-29 / serial 44559) for the same market/product group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 193..256 down to 194..257, opening up a blank row 193.
$ws.Rows.Item(193).Insert()

# Populate the newly opened row 193 with the new record.
$ws.Range("A193").Value = 8
$ws.Range("B193").Value = "Terminal La Palmera de La Serena"
$ws.Range("C193").Value = "Coquimbo"
$ws.Range("D193").Value = 44559
$ws.Range("E193").Value = 4
$ws.Range("F193").Value = 100114013
$ws.Range("G193").Value = "Zanahoria"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 800
$ws.Range("K193").Value = 6000
$ws.Range("L193").Value = 6500
$ws.Range("M193").Value = 6250
$ws.Range("N193").Value = "$/saco 20 kilos"
$ws.Range("O193").Value = "Provincia del Elquí"
$ws.Range("P193").Value = 312
$ws.Range("Q193").Value = 20
$ws.Range("R193").Value = "Hortaliza"
